# Edit script for drimel_only.xlsx style price-list sheet
# Changes:
#  1. D6: 2944 -> 3040
#  2. Insert a new row at row 16 (pushes current rows 16-23 down to 17-24)
#     with the new "PAMPERS PREMIUM XTR RN x36" item.
#  3. Append two brand new rows (25 & 26) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (not auto-converted to a number) while
# keeping the cell's existing "General" number format, so numeric-looking
# strings (SKU codes, barcodes) stay text just like in the source data.
function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

# --- 1. Update existing cost value -----------------------------------
$ws.Range("D6").Value = 3040

# --- 2. Insert new row 16 ---------------------------------------------
$ws.Rows("16:16").Insert()

Set-TextValue $ws.Range("A16") "PAMPERS PREMIUM XTR RN x36"
Set-TextValue $ws.Range("B16") "12972"
Set-TextValue $ws.Range("C16") "7500435132954"
$ws.Range("D16").Value = 1950

# --- 3. Append two new rows at the bottom of the table -----------------
# These rows sit below the previous used range, so they start out with the
# default (non-centered) style. Match the sheet's existing look (centered)
# before populating them.
$ws.Range("A25:D26").HorizontalAlignment = -4108

Set-TextValue $ws.Range("A25") "PAMPERS BabyDry MES XXGx54"
Set-TextValue $ws.Range("B25") "16583"
Set-TextValue $ws.Range("C25") "7500435228725"
$ws.Range("D25").Value = 6608

Set-TextValue $ws.Range("A26") "BABYSEC TOALLAS PREMIUMx50"
Set-TextValue $ws.Range("B26") "16594"
Set-TextValue $ws.Range("C26") "7806500730514"
$ws.Range("D26").Value = 629
